$p = $ppt.ActivePresentation
$s = $p.Slides.Item(14)
$shape = $s.Shapes.Item(2)
$tf = $shape.TextFrame
$tr = $tf.TextRange

$para5 = $tr.Paragraphs(5, 1)
$para6 = $tr.Paragraphs(6, 1)

$para5.Runs(3, 1).Text = " Exercise 05 – IoT Data Processing"
$para6.Runs(3, 1).Text = " Exercise 06 – IoT Security"
